# Update weekly price data for rows 4-15 (Fruta / hortaliza, semanal)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns:
#   D = Fecha (date serial), I = Calidad, J = Volumen,
#   K = Precio mínimo, L = Precio máximo, M = Precio promedio ponderado,
#   P = Precio $/Kg  (mirrors M in this dataset)
$rows = @(
    @{ Row = 4;  D = 44174; I = "Segunda"; J = 800;  K = 450; L = 500; M = 475; P = 475 }
    @{ Row = 5;  D = 44174; I = "Tercera"; J = 1200; K = 250; L = 350; M = 300; P = 300 }
    @{ Row = 6;  D = 44245; I = "Primera"; J = 800;  K = 850; L = 900; M = 875; P = 875 }
    @{ Row = 7;  D = 44245; I = "Segunda"; J = 1000; K = 750; L = 800; M = 775; P = 775 }
    @{ Row = 8;  D = 44201; I = "Segunda"; J = 500;  K = 800; L = 900; M = 850; P = 850 }
    @{ Row = 9;  D = 44210; I = "Segunda"; J = 900;  K = 600; L = 700; M = 650; P = 650 }
    @{ Row = 10; D = 44278; I = "Segunda"; J = 700;  K = 600; L = 700; M = 650; P = 650 }
    @{ Row = 11; D = 44278; I = "Tercera"; J = 400;  K = 500; L = 600; M = 550; P = 550 }
    @{ Row = 12; D = 44267; I = "Tercera"; J = 400;  K = 500; L = 600; M = 550; P = 550 }
    @{ Row = 13; D = 44229; I = "Segunda"; J = 760;  K = 550; L = 600; M = 575; P = 575 }
    @{ Row = 14; D = 44474; I = "Segunda"; J = 200;  K = 600; L = 700; M = 650; P = 650 }
    @{ Row = 15; D = 44224; I = "Segunda"; J = 800;  K = 850; L = 900; M = 875; P = 875 }
)

foreach ($entry in $rows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 4).Value = $entry.D    # D: Fecha
    $ws.Cells.Item($r, 9).Value = $entry.I    # I: Calidad
    $ws.Cells.Item($r, 10).Value = $entry.J   # J: Volumen
    $ws.Cells.Item($r, 11).Value = $entry.K   # K: Precio mínimo
    $ws.Cells.Item($r, 12).Value = $entry.L   # L: Precio máximo
    $ws.Cells.Item($r, 13).Value = $entry.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $entry.P   # P: Precio $/Kg
}
